# Auto-generated edit script applying numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# across several per-class Leve tables, per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 123.51852
$ws.Range("I11").Value = 123.51852
$ws.Range("K11").Value = 123.51852
$ws.Range("M11").Value = 16.48148
# Row 28
$ws.Range("H28").Value = 686.2105
$ws.Range("I28").Value = 702.6
$ws.Range("J28").Value = 624.75
$ws.Range("K28").Value = 702.6
$ws.Range("L28").Value = 624.75
$ws.Range("M28").Value = -217.6
$ws.Range("N28").Value = -1594.75
# Row 86
$ws.Range("H86").Value = 12149.923
$ws.Range("I86").Value = 12743.75
$ws.Range("K86").Value = 12743.75
$ws.Range("M86").Value = -11620.75
# Row 89
$ws.Range("H89").Value = 12149.923
$ws.Range("I89").Value = 12743.75
$ws.Range("K89").Value = 63718.75
$ws.Range("M89").Value = -58102.75
# Row 107
$ws.Range("H107").Value = 945
$ws.Range("I107").Value = 931.7778
$ws.Range("K107").Value = 931.7778
$ws.Range("M107").Value = 988.2222
# Row 132
$ws.Range("H132").Value = 6894.6313
$ws.Range("I132").Value = 7470.4707
$ws.Range("K132").Value = 22411.4121
$ws.Range("M132").Value = -19881.4121
# Row 138
$ws.Range("H138").Value = 4031.2036
$ws.Range("J138").Value = 4298.93
$ws.Range("L138").Value = 12896.79
$ws.Range("N138").Value = -23176.79

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4281.815
$ws.Range("I2").Value = 2994.3
$ws.Range("K2").Value = 2994.3
$ws.Range("M2").Value = -2881.3
# Row 43
$ws.Range("H43").Value = 44000
$ws.Range("J43").Value = 44000
$ws.Range("L43").Value = 44000
$ws.Range("N43").Value = -44626
# Row 116
$ws.Range("H116").Value = 4281.815
$ws.Range("I116").Value = 2994.3
$ws.Range("K116").Value = 2994.3
$ws.Range("M116").Value = -700.3000000000002
# Row 132
$ws.Range("H132").Value = 3232718.8
$ws.Range("I132").Value = 2041.0667
$ws.Range("K132").Value = 6123.2001
$ws.Range("M132").Value = -3593.2001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4281.815
$ws.Range("I3").Value = 2994.3
$ws.Range("K3").Value = 2994.3
$ws.Range("M3").Value = -2880.3
# Row 45
$ws.Range("H45").Value = 24500
$ws.Range("J45").Value = 24500
$ws.Range("L45").Value = 24500
$ws.Range("N45").Value = -26116
# Row 94
$ws.Range("H94").Value = 754.37933
$ws.Range("I94").Value = 859.5454999999999
$ws.Range("K94").Value = 859.5454999999999
$ws.Range("M94").Value = -408.5454999999999
# Row 134
$ws.Range("H134").Value = 7523.1357
$ws.Range("I134").Value = 2627.3333
$ws.Range("K134").Value = 7881.999899999999
$ws.Range("M134").Value = -5346.999899999999
# Row 135
$ws.Range("H135").Value = 84996.5
$ws.Range("J135").Value = 84996.5
$ws.Range("L135").Value = 84996.5
$ws.Range("N135").Value = -95136.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 21164
$ws.Range("J59").Value = 31523
$ws.Range("L59").Value = 31523
$ws.Range("N59").Value = -33813
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0
# Row 107
$ws.Range("H107").Value = 2159.6924
$ws.Range("J107").Value = 2339.2
$ws.Range("L107").Value = 2339.2
$ws.Range("N107").Value = -6179.2
# Row 132
$ws.Range("H132").Value = 32171216
$ws.Range("I132").Value = 2290.4644
$ws.Range("K132").Value = 6871.3932
$ws.Range("M132").Value = -4341.3932
# Row 138
$ws.Range("H138").Value = 80333.336
$ws.Range("J138").Value = 80333.336
$ws.Range("L138").Value = 80333.336
$ws.Range("N138").Value = -90613.336

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2728805.8
$ws.Range("I4").Value = 3072273.2
$ws.Range("J4").Value = 668000
$ws.Range("K4").Value = 9216819.600000001
$ws.Range("L4").Value = 2004000
$ws.Range("M4").Value = -9216707.600000001
$ws.Range("N4").Value = -2004224
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
# Row 137
$ws.Range("H137").Value = 5111.909
$ws.Range("I137").Value = 3771
$ws.Range("K137").Value = 11313
$ws.Range("M137").Value = -6213

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 19249.25
$ws.Range("I18").Value = 19249.25
$ws.Range("K18").Value = 19249.25
$ws.Range("M18").Value = -18956.25
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = 0
# Row 113
$ws.Range("H113").Value = 969.8570999999999
$ws.Range("I113").Value = 959
$ws.Range("J113").Value = 997
$ws.Range("K113").Value = 959
$ws.Range("L113").Value = 997
$ws.Range("M113").Value = 1211
$ws.Range("N113").Value = -5337
# Row 117
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884
# Row 126
$ws.Range("H126").Value = 5653602
$ws.Range("I126").Value = 2766289.8
$ws.Range("K126").Value = 8298869.399999999
$ws.Range("M126").Value = -8296399.399999999
# Row 132
$ws.Range("I132").Value = 10907.4
$ws.Range("K132").Value = 32722.2
$ws.Range("M132").Value = -30192.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4828
# Row 22
$ws.Range("H22").Value = 62501310
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 142858200
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 142858200
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -142858790
# Row 24
$ws.Range("H24").Value = 2645
$ws.Range("I24").Value = 193.33333
$ws.Range("K24").Value = 193.33333
$ws.Range("M24").Value = 149.66667
# Row 27
$ws.Range("H27").Value = 62501310
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 142858200
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 142858200
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -142858414
# Row 40
$ws.Range("H40").Value = 5348949
$ws.Range("I40").Value = 1491
$ws.Range("K40").Value = 1491
$ws.Range("M40").Value = -1355
# Row 96
$ws.Range("H96").Value = 16666.666
$ws.Range("J96").Value = 16666.666
$ws.Range("L96").Value = 16666.666
$ws.Range("N96").Value = -22158.666
# Row 100
$ws.Range("H100").Value = 2866.375
$ws.Range("I100").Value = 2548.9
$ws.Range("J100").Value = 3395.5
$ws.Range("K100").Value = 2548.9
$ws.Range("L100").Value = 3395.5
$ws.Range("M100").Value = -2007.9
$ws.Range("N100").Value = -4477.5
# Row 132
$ws.Range("H132").Value = 4991295
$ws.Range("I132").Value = 7003
$ws.Range("J132").Value = 6985011.5
$ws.Range("K132").Value = 21009
$ws.Range("L132").Value = 20955034.5
$ws.Range("M132").Value = -18479
$ws.Range("N132").Value = -20960094.5
# Row 139
$ws.Range("H139").Value = 94960
$ws.Range("J139").Value = 94960
$ws.Range("L139").Value = 94960
$ws.Range("N139").Value = -105240

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 647.625
$ws.Range("I107").Value = 689.65
$ws.Range("K107").Value = 2068.95
$ws.Range("M107").Value = -148.9499999999998
# Row 113
$ws.Range("H113").Value = 3469.476
$ws.Range("I113").Value = 3486.5293
$ws.Range("K113").Value = 10459.5879
$ws.Range("M113").Value = -8289.5879
# Row 132
$ws.Range("H132").Value = 746606.9
$ws.Range("I132").Value = 11053.9
$ws.Range("K132").Value = 33161.7
$ws.Range("M132").Value = -30631.7
